$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.481.48"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.828.04"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'315.11"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -4.42%  "
$ws.Range("D8").Value = "'0.3908"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "'0.07648"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'41.85"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'1.109"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "'21.10"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "'7.585"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "1.824.58"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "'93.13"
$ws.Range("E17").Value = "  +5.51%  "
$ws.Range("D18").Value = "'0.00001083"
$ws.Range("D19").Value = "'0.06666"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").Value = "'17.73"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'6.165"
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("D23").Value = "28.505.95"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +7.67%  "
$ws.Range("D26").Value = "'156.93"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'20.64"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").Value = "2.034.50"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "'2.395"
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("D30").Value = "'125.15"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "'1.127"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").Value = "'0.1085"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.689"
$ws.Range("E33").Value = "  +3.84%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'0.07042"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'0.2226"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").Value = "'8.970"
$ws.Range("E37").Value = "  +6.88%  "
$ws.Range("D38").Value = "'0.02323"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "'5.141"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "'0.6274"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'1.396"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'13.39"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'0.5911"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "'3.715"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'124.47"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'1.984"
$ws.Range("E49").Value = "  +3.82%  "
$ws.Range("D50").Value = "'1.194"
$ws.Range("D51").Value = "'0.06921"
$ws.Range("E51").Value = "  +1.69%  "
